# Apply "Updated data to reflect new requirement separation" to the RPTA sheet.
# This inserts three new columns (Corequisites, Concurrent, Recommended) between
# the existing "Prerequisites" (C) and "Terms Typically Offered" (D) columns,
# populates them, and cleans up a couple of "Recommended:" notes that had been
# embedded inside the Prerequisites text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert three new columns before column D. This shifts the existing
#    "Terms Typically Offered" column from D to G automatically.
$ws.Columns("D:F").Insert()

# 2. Set the new header row values.
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# 3. Default-fill the new Corequisites/Concurrent/Recommended columns with "NA"
#    for every data row (2-56).
$ws.Range("D2:D56").Value = "NA"
$ws.Range("E2:E56").Value = "NA"
$ws.Range("F2:F56").Value = "NA"

# 4. A couple of rows had a "Recommended: ..." note tacked on to the end of the
#    Prerequisites text. Pull that text out into the new Recommended column and
#    trim it off the Prerequisites column. (The Terms value picks up a trailing
#    space left behind from the original formatting.)
$ws.Range("C40").Value = "AGB 214 or BUS 212; AGB 323 or BUS 215; RPTA 360 with C- or better or consent of instructor; and senior standing."
$ws.Range("F40").Value = "ENGL 310."
$ws.Range("G40").Value = "F, W, SP "

$ws.Range("C44").Value = "RPTA majors only and senior standing."
$ws.Range("F44").Value = "enrollment two quarters prior to RPTA 465."
$ws.Range("G44").Value = "F, W, SP "

# 5. Clean up wording in the Prerequisites column for RPTA 323.
$ws.Range("C25").Value = "one of the RPTA major, Event Planning and Experience Management minor, or Exercise and Sport Studies minor; and junior standing."

# 6. Normalize stray non-breaking spaces within course-number references (e.g.
#    "RPTA 210") in the Prerequisites column to regular spaces.
$nbspRows = @(18,19,20,24,26,28,29,31,33,34,35,36,37,38,39,42,43)
foreach ($r in $nbspRows) {
    $cell = $ws.Range("C$r")
    $text = $cell.Text
    $cell.Value = $text.Replace([char]160, ' ')
}
